# Weekly update: insert a new weekly record as row 4 and push the
# previously existing rows (old rows 4-21) down by one (new rows 5-22).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 4; this shifts old rows 4..21 -> 5..22
$ws.Rows.Item(4).Insert()

# Fill the constant (record-identity) columns for the new row, copied
# from the surrounding records (row 5, the former row 4).
$ws.Range("A4").Value = 7
$ws.Range("B4").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C4").Value = "Ñuble"
$ws.Range("D4").Value = 45099
$ws.Range("E4").Value = 16
$ws.Range("F4").Value = "Fruta"
$ws.Range("G4").Value = 100107
$ws.Range("H4").Value = "Otros"
$ws.Range("I4").Value = 100107011
$ws.Range("J4").Value = "Tuna"
$ws.Range("K4").Value = "Sin especificar"
$ws.Range("L4").Value = "Primera"
$ws.Range("M4").Value = 40
$ws.Range("N4").Value = 22000
$ws.Range("O4").Value = 22000
$ws.Range("P4").Value = 22000
$ws.Range("Q4").Value = "`$/caja 18 kilos"
$ws.Range("R4").Value = "Región Metropolitana"
$ws.Range("S4").Value = 1222
$ws.Range("T4").Value = 18
